$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Ludzie" to "humanAncestry"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "humanAncestry"

# Update the selection on sheet1 to D35
$ws1.Activate()
$ws1.Range("D35").Select()
